$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" positioned right before the "总计"
#    sheet (i.e. after "2021-Q4"), matching the target sheet order:
#    2020-Q4, 2021-Q4, 2022-Q1, 总计
# ---------------------------------------------------------------------------
$beforeRef = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($beforeRef)
$q1Sheet.Name = "2022-Q1"

# NOTE: after Add()+rename, re-fetch sheets fresh by name rather than
# reusing older variable references, which can become stale/aliased.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header styling (bold font + border) from the "2021-Q4" sheet header
# row so the new sheet matches the look of its siblings.
$srcHeader = $wb.Worksheets.Item("2021-Q4").Range("B1:H1")
$srcHeader.Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

# The fund-code / numeric-looking text columns (B-G) must stay text so
# leading zeros ("006478") and trailing zeros ("0.0901") are preserved
# exactly like the source data (stored as inline strings, not numbers).
$q1Sheet.Range("B1:G3").NumberFormat = "@"

# Copy the index-column style (bold font + border) used on A2/A3 from an
# existing sheet onto the new sheet's A2/A3 cells.
$srcIndex = $wb.Worksheets.Item("2021-Q4").Range("A2:A3")
$srcIndex.Copy()
$q1Sheet.Range("A2:A3").PasteSpecial(-4122)

# Header row
$q1Sheet.Cells.Item(1, 2).Value = "基金代码"
$q1Sheet.Cells.Item(1, 3).Value = "基金名称"
$q1Sheet.Cells.Item(1, 4).Value = "基金规模"
$q1Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q1Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q1Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1Sheet.Cells.Item(1, 8).Value = "仓位排名"

# Row 2 data - 510081 长盛动态精选混合
$q1Sheet.Cells.Item(2, 1).Value = 0
$q1Sheet.Cells.Item(2, 2).Value = "510081"
$q1Sheet.Cells.Item(2, 3).Value = "长盛动态精选混合"
$q1Sheet.Cells.Item(2, 4).Value = "3.15"
$q1Sheet.Cells.Item(2, 5).Value = "60.76"
$q1Sheet.Cells.Item(2, 6).Value = "2.86"
$q1Sheet.Cells.Item(2, 7).Value = "0.0901"
$q1Sheet.Cells.Item(2, 8).Value = 9

# Row 3 data - 006478 长盛多因子策略优选股票
$q1Sheet.Cells.Item(3, 1).Value = 1
$q1Sheet.Cells.Item(3, 2).Value = "006478"
$q1Sheet.Cells.Item(3, 3).Value = "长盛多因子策略优选股票"
$q1Sheet.Cells.Item(3, 4).Value = "0.51"
$q1Sheet.Cells.Item(3, 5).Value = "84.41"
$q1Sheet.Cells.Item(3, 6).Value = "4.33"
$q1Sheet.Cells.Item(3, 7).Value = "0.0221"
$q1Sheet.Cells.Item(3, 8).Value = 7

# ---------------------------------------------------------------------------
# 2. Add a new top data row to the "总计" sheet summarising the 2022-Q1
#    sheet (2 funds held, 0.11 亿元 total market value), pushing the
#    existing rows down and renumbering the index column.
# ---------------------------------------------------------------------------
$totalSheet.Rows(2).Insert()

# Re-apply the bold/border style to the new index cell (A2); the plain
# data cells (B2:D2) should stay unstyled like the rest of the data rows.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.11

# Renumber the index column for the rows that shifted down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
